# Update policies data download file: append the data-source citation
# block below the existing table and restore the header row's wrap height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 is left blank (matches the original author's layout), then three
# new rows are appended with source information for the dataset.

# Write the URL first so it lands in the shared-string table ahead of the
# "Source: " label and description strings, matching the author's order.
$ws.Range("A23").Value = "http://www.edr.state.fl.us/Content/conferences/criminaljustice/trends.pdf"

$ws.Range("A21").Value = "Source: "
$ws.Range("A21").Font.Bold = $true

$ws.Range("A22").Value = "Florida Office of Economic and Demographic Research: Criminal Justice Trends (February 23, 2017)"

# The header row (row 2) wraps text; restore its height now that the sheet
# has been touched again.
$ws.Rows("2").RowHeight = 45

# Leave the selection where the author left it when they saved the file.
$ws.Range("F28").Select()
